# Apply text replacements for the updated worksheet (date + division problems)

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-03-13 Wednesday"; new = "2024-03-14 Thursday"},
    @{old = "703÷4="; new = "271÷7="},
    @{old = "579÷8="; new = "796÷7="},
    @{old = "248÷2="; new = "112÷4="},
    @{old = "157÷9="; new = "922÷3="},
    @{old = "119÷5="; new = "609÷5="},
    @{old = "795÷3="; new = "778÷8="},
    @{old = "951÷7="; new = "231÷7="},
    @{old = "571÷8="; new = "814÷7="},
    @{old = "795÷2="; new = "182÷3="},
    @{old = "180÷9="; new = "440÷6="},
    @{old = "688÷8="; new = "734÷6="},
    @{old = "393÷9="; new = "232÷7="},
    @{old = "303÷9="; new = "110÷5="},
    @{old = "538÷6="; new = "764÷4="},
    @{old = "226÷5="; new = "419÷7="},
    @{old = "782÷2="; new = "128÷9="},
    @{old = "410÷8="; new = "904÷4="},
    @{old = "760÷5="; new = "587÷7="},
    @{old = "417÷6="; new = "102÷8="},
    @{old = "609÷6="; new = "327÷2="},
    @{old = "166÷8="; new = "466÷6="},
    @{old = "782÷4="; new = "292÷9="},
    @{old = "110÷7="; new = "314÷2="},
    @{old = "549÷4="; new = "104÷9="},
    @{old = "712÷5="; new = "973÷6="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
